# Fruta / hortaliza, semanal
#
# Insert one new weekly price observation as a new row 350 on the
# "Pomelo" sheet (Vega Modelo de Temuco / Start Ruby / Primera),
# pushing the existing rows 350-362 down to 351-363.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 350 - shifts rows 350..362 down to 351..363
# and extends the used range / dimension to A1:T363 automatically.
$ws.Rows(350).Insert()

# Populate the newly inserted row 350 with the new observation.
$ws.Cells.Item(350, 1).Value  = 10
$ws.Cells.Item(350, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(350, 3).Value  = "La Araucanía"
$ws.Cells.Item(350, 4).Value  = 44939
$ws.Cells.Item(350, 5).Value  = 9
$ws.Cells.Item(350, 6).Value  = "Fruta"
$ws.Cells.Item(350, 7).Value  = 100102
$ws.Cells.Item(350, 8).Value  = "Cítricos"
$ws.Cells.Item(350, 9).Value  = 100102006
$ws.Cells.Item(350, 10).Value = "Pomelo"
$ws.Cells.Item(350, 11).Value = "Start Ruby"
$ws.Cells.Item(350, 12).Value = "Primera"
$ws.Cells.Item(350, 13).Value = 100
$ws.Cells.Item(350, 14).Value = 15000
$ws.Cells.Item(350, 15).Value = 15000
$ws.Cells.Item(350, 16).Value = 15000
$ws.Cells.Item(350, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(350, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(350, 19).Value = 1000
$ws.Cells.Item(350, 20).Value = 15
